# wmt_etl/tests/data/full_inputs/ND02.xlsx update:
# "Update XLSXs to have no blank columns for now"

$wb = $excel.ActiveWorkbook

$wsExtract = $wb.Worksheets.Item("WMT_Extract")
$wsInst    = $wb.Worksheets.Item("Inst_Reports")

# --- sharedStrings.xml: "ND02" -> "ND01" (used by WMT_Extract!C2 and C3) ---
$wsExtract.Range("C2").Value = "ND01"
$wsExtract.Range("C3").Value = "ND01"

# --- Inst_Reports!G2: fill in the previously-blank column (matches row 3's F/G pair) ---
$wsInst.Range("G2").Value = 2

# --- View state ---
# WMT_Extract: no longer the selected tab; scrolled/selected cell moves from AE1 to AE2
# Inst_Reports: becomes the selected tab; selection moves from D26 to E14
$wsExtract.Activate()
$wsExtract.Range("AE2").Select()

$wsInst.Activate()
$wsInst.Range("E14").Select()
